# Update countries & provincias Spain
# Applies the 21-Mar-2020 17:41 data refresh to the "Pais" sheet:
#  - Updates the "Datos actualizados" timestamp in A1
#  - Austria now ranks above Belgica (new Austria numbers; Belgica keeps its old numbers)
#  - Argelia now ranks above Letonia (new Argelia numbers; Letonia/Costa Rica/
#    Republica Dominicana/Uruguay/Hungria each shift down one row, keeping
#    their previous numbers)
#  - Several other countries (Estados Unidos, Suiza, Reino Unido, Brasil,
#    Japon, Chequia, India) got refreshed case/death counts

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp header
$ws.Range("A1").Value = "Datos actualizados a 21 de Marzo de 2020 a las 17:41"

# Helper: write a full data row (Pais, Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes)
function Set-CountryRow($row, $pais, $casosTotales, $nuevosCasos, $casosActivos, $recuperados, $casosCriticos, $muertesHoy, $muertes) {
    $ws.Range("A$row").Value = $pais
    $ws.Range("B$row").Value = $casosTotales
    $ws.Range("C$row").Value = $nuevosCasos
    $ws.Range("D$row").Value = $casosActivos
    $ws.Range("E$row").Value = $recuperados
    $ws.Range("F$row").Value = $casosCriticos
    $ws.Range("G$row").Value = $muertesHoy
    $ws.Range("H$row").Value = $muertes
}

# Row 7: Estados Unidos - refreshed counts
Set-CountryRow 7 "Estados Unidos" 22085 2702 147 21656 64 26 282

# Row 12: Suiza - refreshed counts
Set-CountryRow 12 "Suiza" 6192 577 15 6113 141 8 64

# Row 13: Reino Unido - refreshed counts
Set-CountryRow 13 "Reino Unido" 4094 111 65 3849 20 3 180

# Row 15-16: Austria now overtakes Belgica in the ranking.
# Austria gets new, updated numbers; Belgica keeps its previous numbers.
Set-CountryRow 15 "Austria" 2847 198 9 2830 15 2 8
Set-CountryRow 16 "Belgica" 2815 558 263 2485 288 30 67

# Row 17: Noruega - refreshed counts
Set-CountryRow 17 "Noruega" 2082 123 1 2074 28 0 7

# Row 24: Brasil - refreshed counts
Set-CountryRow 24 "Brasil" 1021 51 2 1001 18 7 18

# Row 25: Japon - refreshed counts
Set-CountryRow 25 "Japon" 1007 0 215 757 55 0 35

# Row 26: Chequia - refreshed counts
Set-CountryRow 26 "Chequia" 925 92 6 919 7 0 0

# Row 46: India - refreshed counts
Set-CountryRow 46 "India" 329 80 23 301 0 0 5

# Rows 69-75: Argelia now overtakes Letonia, Costa Rica, Republica Dominicana,
# Uruguay and Hungria in the ranking. Eslovaquia and Argelia get new, updated
# numbers; the others keep their previous numbers, shifted down one row.
Set-CountryRow 69 "Eslovaquia" 146 9 7 139 2 0 0
Set-CountryRow 70 "Argelia" 139 45 65 59 0 4 15
Set-CountryRow 71 "Letonia" 124 13 1 123 0 0 0
Set-CountryRow 72 "Costa Rica" 113 0 2 109 2 0 2
Set-CountryRow 73 "Republica Dominicana" 112 40 0 109 0 1 3
Set-CountryRow 74 "Uruguay" 110 0 0 110 0 0 0
Set-CountryRow 75 "Hungria" 103 18 7 92 6 0 4
